$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values that look like
# numbers (e.g. "536.82", "0.998") are stored as text, matching the
# original inlineStr cell type rather than being auto-coerced to Number.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.365.69'
$ws.Range("E2").Value = '  -5.67%  '
$ws.Range("D3").Value = '2.456.69'
$ws.Range("E3").Value = '  -8.33%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '536.82'
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("D6").Value = '148.12'
$ws.Range("E6").Value = '  -6.56%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -4.08%  '
$ws.Range("D9").Value = '2.473.63'
$ws.Range("E9").Value = '  -7.79%  '
$ws.Range("E10").Value = '  -5.73%  '
$ws.Range("E11").Value = '  -2.49%  '
$ws.Range("D12").Value = '5.35'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("E13").Value = '  -4.09%  '
$ws.Range("D14").Value = '2.905.21'
$ws.Range("E14").Value = '  -7.85%  '
$ws.Range("D15").Value = '24.12'
$ws.Range("E15").Value = '  -8.24%  '
$ws.Range("D16").Value = '59.350.98'
$ws.Range("E16").Value = '  -5.51%  '
$ws.Range("E17").Value = '  -5.82%  '
$ws.Range("D18").Value = '2.495.77'
$ws.Range("E18").Value = '  -6.84%  '
$ws.Range("D19").Value = '11.18'
$ws.Range("E19").Value = '  -5.82%  '
$ws.Range("E20").Value = '  -5.44%  '
$ws.Range("D21").Value = '324.16'
$ws.Range("E21").Value = '  -6.15%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '0.964'
$ws.Range("E22").Value = '  -3.53%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  -8.61%  '
$ws.Range("D24").Value = '0.463'
$ws.Range("E24").Value = '  -8.57%  '
$ws.Range("D25").Value = '60.64'
$ws.Range("E25").Value = '  -4.26%  '
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").Value = '0.977'
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("D28").Value = '7.74'
$ws.Range("E28").Value = '  -5.57%  '
$ws.Range("D29").Value = '6.82'
$ws.Range("E29").Value = '  -5.52%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '1.27'
$ws.Range("E30").Value = '  -8.30%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '1.82'
$ws.Range("E31").Value = '  -6.08%  '
$ws.Range("E32").Value = '  -9.92%  '
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '155.49'
$ws.Range("E34").Value = '  -6.44%  '
$ws.Range("D35").Value = '1.39'
$ws.Range("E35").Value = '  -5.68%  '
$ws.Range("D36").Value = '4.56'
$ws.Range("E36").Value = '  -6.19%  '
$ws.Range("E37").Value = '  -5.68%  '
$ws.Range("E38").Value = '  -1.52%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '5.93'
$ws.Range("E39").Value = '  -6.37%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '318.74'
$ws.Range("E40").Value = '  -9.35%  '
$ws.Range("D41").Value = '36.75'
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("D42").Value = '0.843'
$ws.Range("E42").Value = '  -12.86%  '
$ws.Range("E43").Value = '  -7.09%  '
$ws.Range("D44").Value = '0.995'
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '10.74'
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("E46").Value = '  -5.18%  '
$ws.Range("E47").Value = '  -3.33%  '
$ws.Range("D48").Value = '0.0527'
$ws.Range("E48").Value = '  -6.36%  '
$ws.Range("D49").Value = '19.08'
$ws.Range("E49").Value = '  -8.96%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '18.54'
$ws.Range("E50").Value = '  -8.91%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '0.0229'
$ws.Range("E51").Value = '  -5.43%  '
